$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.877.49'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.293.53'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '98.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E7').Value = '  -1.20%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.606'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.84'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '2.636.22'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.857'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '2.295.48'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '43.778.66'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('E20').Value = '  -3.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('E22').Value = '  +7.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.15%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('E32').Value = '  -3.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0893'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.44'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.90%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('E39').Value = '  +4.21%  '
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.54%  '
$ws.Range('E45').Value = '  -3.51%  '
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.442'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.11%  '
$ws.Range('E51').Value = '  +11.40%  '
